$d = $word.ActiveDocument

# Update the header date range text (start/end dates)
# and all 31 daily date cells in the attendance table.
# Each old date maps 1:1 to a new date, 380 days later -> +1 year +10 days.
$d.Content.Find.Execute("2019-07-26", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-05", 2) | Out-Null
$d.Content.Find.Execute("2019-07-27", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-06", 2) | Out-Null
$d.Content.Find.Execute("2019-07-28", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-07", 2) | Out-Null
$d.Content.Find.Execute("2019-07-29", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-08", 2) | Out-Null
$d.Content.Find.Execute("2019-07-30", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-09", 2) | Out-Null
$d.Content.Find.Execute("2019-07-31", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-10", 2) | Out-Null
$d.Content.Find.Execute("2019-08-01", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-11", 2) | Out-Null
$d.Content.Find.Execute("2019-08-02", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-12", 2) | Out-Null
$d.Content.Find.Execute("2019-08-03", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-13", 2) | Out-Null
$d.Content.Find.Execute("2019-08-04", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-14", 2) | Out-Null
$d.Content.Find.Execute("2019-08-05", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-15", 2) | Out-Null
$d.Content.Find.Execute("2019-08-06", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-16", 2) | Out-Null
$d.Content.Find.Execute("2019-08-07", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-17", 2) | Out-Null
$d.Content.Find.Execute("2019-08-08", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-18", 2) | Out-Null
$d.Content.Find.Execute("2019-08-09", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-19", 2) | Out-Null
$d.Content.Find.Execute("2019-08-10", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-20", 2) | Out-Null
$d.Content.Find.Execute("2019-08-11", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-21", 2) | Out-Null
$d.Content.Find.Execute("2019-08-12", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-22", 2) | Out-Null
$d.Content.Find.Execute("2019-08-13", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-23", 2) | Out-Null
$d.Content.Find.Execute("2019-08-14", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-24", 2) | Out-Null
$d.Content.Find.Execute("2019-08-15", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-25", 2) | Out-Null
$d.Content.Find.Execute("2019-08-16", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-26", 2) | Out-Null
$d.Content.Find.Execute("2019-08-17", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-27", 2) | Out-Null
$d.Content.Find.Execute("2019-08-18", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-28", 2) | Out-Null
$d.Content.Find.Execute("2019-08-19", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-29", 2) | Out-Null
$d.Content.Find.Execute("2019-08-20", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-30", 2) | Out-Null
$d.Content.Find.Execute("2019-08-21", $true, $false, $false, $false, $false, $true, 1, $false, "2020-08-31", 2) | Out-Null
$d.Content.Find.Execute("2019-08-22", $true, $false, $false, $false, $false, $true, 1, $false, "2020-09-01", 2) | Out-Null
$d.Content.Find.Execute("2019-08-23", $true, $false, $false, $false, $false, $true, 1, $false, "2020-09-02", 2) | Out-Null
$d.Content.Find.Execute("2019-08-24", $true, $false, $false, $false, $false, $true, 1, $false, "2020-09-03", 2) | Out-Null
$d.Content.Find.Execute("2019-08-25", $true, $false, $false, $false, $false, $true, 1, $false, "2020-09-04", 2) | Out-Null

# Update every row height in the attendance table from 369 -> 397 twips
# (18.45pt -> 19.85pt).
$table = $d.Tables(1)
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $table.Rows($i).Height = 19.85
}
